$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 488.66666
$ws.Range("I12").Value = 365
$ws.Range("J12").Value = 550.5
$ws.Range("K12").Value = 365
$ws.Range("L12").Value = 550.5
$ws.Range("M12").Value = -195
$ws.Range("N12").Value = -890.5

$ws.Range("H55").Value = 373.5
$ws.Range("I55").Value = 373.5
$ws.Range("K55").Value = 373.5
$ws.Range("M55").Value = -159.5

$ws.Range("H88").Value = 3608.6
$ws.Range("J88").Value = 4157.2
$ws.Range("L88").Value = 4157.2
$ws.Range("N88").Value = -4969.2

$ws.Range("H91").Value = 3608.6
$ws.Range("J91").Value = 4157.2
$ws.Range("L91").Value = 4157.2
$ws.Range("N91").Value = -6965.2

$ws.Range("H129").Value = 3031.5715
$ws.Range("I129").Value = 2544.2
$ws.Range("J129").Value = 4250
$ws.Range("K129").Value = 7632.599999999999
$ws.Range("L129").Value = 12750
$ws.Range("M129").Value = -2632.599999999999
$ws.Range("N129").Value = -22750

$ws.Range("H131").Value = 816
$ws.Range("I131").Value = 690
$ws.Range("J131").Value = 1194
$ws.Range("K131").Value = 2070
$ws.Range("L131").Value = 3582
$ws.Range("M131").Value = 2970
$ws.Range("N131").Value = -13662

$ws.Range("H138").Value = 1238
$ws.Range("I138").Value = 1202.3529
$ws.Range("J138").Value = 1281.2858
$ws.Range("K138").Value = 3607.0587
$ws.Range("L138").Value = 3843.8574
$ws.Range("M138").Value = 1532.9413
$ws.Range("N138").Value = -14123.8574

$ws.Range("H139").Value = 94000
$ws.Range("I139").Value = 94000
$ws.Range("K139").Value = 94000
$ws.Range("M139").Value = -88860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 462.75
$ws.Range("I2").Value = 423
$ws.Range("K2").Value = 423
$ws.Range("M2").Value = -310

$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2000
$ws.Range("K33").Value = 2000
$ws.Range("M33").Value = -1671

$ws.Range("H45").Value = 2533
$ws.Range("I45").Value = 2388.3333
$ws.Range("K45").Value = 2388.3333
$ws.Range("M45").Value = -2011.3333

$ws.Range("H74").Value = 3356.8572
$ws.Range("I74").Value = 2546.8572
$ws.Range("K74").Value = 2546.8572
$ws.Range("M74").Value = -1672.8572

$ws.Range("H77").Value = 3356.8572
$ws.Range("I77").Value = 2546.8572
$ws.Range("K77").Value = 12734.286
$ws.Range("M77").Value = -8366.286

$ws.Range("H97").Value = 2111
$ws.Range("I97").Value = 1016.6
$ws.Range("K97").Value = 1016.6
$ws.Range("M97").Value = -520.6

$ws.Range("H116").Value = 462.75
$ws.Range("I116").Value = 423
$ws.Range("K116").Value = 423
$ws.Range("M116").Value = 1871

$ws.Range("H132").Value = 4077.6667
$ws.Range("I132").Value = 3374.5
$ws.Range("K132").Value = 10123.5
$ws.Range("M132").Value = -7593.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 462.75
$ws.Range("I3").Value = 423
$ws.Range("K3").Value = 423
$ws.Range("M3").Value = -309

$ws.Range("H86").Value = 10828.6
$ws.Range("J86").Value = 13329.5
$ws.Range("L86").Value = 13329.5
$ws.Range("N86").Value = -15575.5

$ws.Range("H89").Value = 10828.6
$ws.Range("J89").Value = 13329.5
$ws.Range("L89").Value = 66647.5
$ws.Range("N89").Value = -77879.5

$ws.Range("H99").Value = 1380.375
$ws.Range("I99").Value = 1020.5714
$ws.Range("K99").Value = 1020.5714
$ws.Range("M99").Value = 477.4286

$ws.Range("H134").Value = 350
$ws.Range("I134").Value = 350
$ws.Range("K134").Value = 1050
$ws.Range("M134").Value = 1485

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4433.3335
$ws.Range("I58").Value = 1650
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 1650
$ws.Range("L58").Value = 10000
$ws.Range("M58").Value = -1447
$ws.Range("N58").Value = -10406

$ws.Range("H62").Value = 3448.6
$ws.Range("I62").Value = 3310.75
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3310.75
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2686.75
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 3448.6
$ws.Range("I65").Value = 3310.75
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 16553.75
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -13433.75
$ws.Range("N65").Value = -26240

$ws.Range("H136").Value = 4433.3335
$ws.Range("I136").Value = 1650
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 4950
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -2400
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 74503
$ws.Range("I7").Value = 125159.125
$ws.Range("J7").Value = 6961.5
$ws.Range("K7").Value = 375477.375
$ws.Range("L7").Value = 20884.5
$ws.Range("M7").Value = -375365.375
$ws.Range("N7").Value = -21108.5

$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 12000
$ws.Range("M86").Value = -10814

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H88").Value = 19999
$ws.Range("J88").Value = 19999
$ws.Range("L88").Value = 59997
$ws.Range("N88").Value = -60853

$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 36000
$ws.Range("M89").Value = -30072

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H91").Value = 19999
$ws.Range("J91").Value = 19999
$ws.Range("L91").Value = 59997
$ws.Range("N91").Value = -62961

$ws.Range("H122").Value = 190
$ws.Range("J122").Value = 190
$ws.Range("L122").Value = 1710
$ws.Range("N122").Value = -6610

$ws.Range("H124").Value = 1000
$ws.Range("I124").Value = 1000
$ws.Range("K124").Value = 3000
$ws.Range("M124").Value = 1910

$ws.Range("H126").Value = 6350
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = 1340

$ws.Range("H129").Value = 1599.5
$ws.Range("J129").Value = 1200
$ws.Range("L129").Value = 3600
$ws.Range("N129").Value = -13600

$ws.Range("H131").Value = 2437.375
$ws.Range("I131").Value = 1666.6666
$ws.Range("J131").Value = 2899.8
$ws.Range("K131").Value = 4999.9998
$ws.Range("L131").Value = 8699.400000000001
$ws.Range("M131").Value = 40.0002000000004
$ws.Range("N131").Value = -18779.4

$ws.Range("H132").Value = 5000
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 45000
$ws.Range("N132").Value = -50060

$ws.Range("H134").Value = 1011
$ws.Range("I134").Value = 1011
$ws.Range("K134").Value = 3033
$ws.Range("M134").Value = 2037

$ws.Range("H139").Value = 7816.5
$ws.Range("I139").Value = 800
$ws.Range("K139").Value = 2400
$ws.Range("M139").Value = 2740

$ws.Range("H140").Value = 1429855.6
$ws.Range("I140").Value = 1429855.6
$ws.Range("K140").Value = 4289566.800000001
$ws.Range("M140").Value = -4284386.800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 921.7143
$ws.Range("I102").Value = 642
$ws.Range("K102").Value = 642
$ws.Range("M102").Value = 980

$ws.Range("H132").Value = 2403.36
$ws.Range("I132").Value = 2397.3
$ws.Range("K132").Value = 7191.900000000001
$ws.Range("M132").Value = -4661.900000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1046.579
$ws.Range("I55").Value = 319.16666
$ws.Range("J55").Value = 1382.3077
$ws.Range("K55").Value = 319.16666
$ws.Range("L55").Value = 1382.3077
$ws.Range("M55").Value = -146.16666
$ws.Range("N55").Value = -1728.3077

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 55000
$ws.Range("J135").Value = 55000
$ws.Range("L135").Value = 55000
$ws.Range("N135").Value = -65140

$ws.Range("H136").Value = 23960.5
$ws.Range("I136").Value = 1753.375
$ws.Range("K136").Value = 5260.125
$ws.Range("M136").Value = -2710.125

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2753.8948
$ws.Range("I132").Value = 2195.5293
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 6586.5879
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -4056.5879
$ws.Range("N132").Value = -27560
